# Added constraints to visit forms
#
# The "survey" sheet gains two new columns: "constraint" (H) and
# "constraint_message" (I). Two existing question rows (plant_height and
# crop_yield) get a validation expression + a human readable error message.

$wb = $excel.ActiveWorkbook

# --- queries sheet: restore the last-used selection (A2) -----------------
$queries = $wb.Worksheets.Item("queries")
$queries.Cells.Item(2, 1).Select() | Out-Null

# --- survey sheet: add the constraint / constraint_message columns -------
$survey = $wb.Worksheets.Item("survey")

# Header row
$survey.Cells.Item(1, 8).Value = "constraint"
$survey.Cells.Item(1, 9).Value = "constraint_message"

# crop_yield (row 34) and plant_height (row 4) messages are authored first,
# then the corresponding constraint formulas - this keeps the shared-string
# table ordering identical to the authored workbook.
$survey.Cells.Item(34, 9).Value = "The crop yield must be a positive value less than 8000 kg. Please enter a valid number."
$survey.Cells.Item(4, 9).Value = "The height of a maize plant must be a positive value less than 1100 cm. Please enter a valid number."
$survey.Cells.Item(4, 8).Value = "data('plant_height')  >  0 && data('plant_height') < 1100"
$survey.Cells.Item(34, 8).Value = "data('crop_yield') > 0 && data('crop_yield')  <  8000"

# Widen column G (display.text) and size the two new columns.
$survey.Columns.Item(7).ColumnWidth = 58.25
$survey.Columns.Item(8).ColumnWidth = 18.25
$survey.Columns.Item(9).ColumnWidth = 19.45

# Leave the cursor on the last-edited cell and make "survey" the active tab.
$survey.Cells.Item(34, 8).Select() | Out-Null
